$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 479, shifting the existing 479:560 data block
# down to 482:563 (dimension grows from A1:R560 to A1:R563).
$ws.Range("A479:A481").EntireRow.Insert()

# Populate the 3 newly inserted rows with the new week's price report
# (Fecha = 44505), following the same Extra / Primera / Segunda layout
# used throughout the rest of the sheet.
$newRows = @(
    @(479, "Extra",    130, 13000, 13000, 13000, 4333),
    @(480, "Primera",  150, 10000, 10000, 10000, 3333),
    @(481, "Segunda",  100,  8000,  8000,  8000, 2667)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value  = 6
    $ws.Cells.Item($rowNum, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($rowNum, 3).Value  = "Metropolitana"
    $ws.Cells.Item($rowNum, 4).Value  = 44505
    $ws.Cells.Item($rowNum, 5).Value  = 13
    $ws.Cells.Item($rowNum, 6).Value  = 100112009
    $ws.Cells.Item($rowNum, 7).Value  = "Acelga"
    $ws.Cells.Item($rowNum, 8).Value  = "Sin especificar"
    $ws.Cells.Item($rowNum, 9).Value  = $r[1]
    $ws.Cells.Item($rowNum, 10).Value = $r[2]
    $ws.Cells.Item($rowNum, 11).Value = $r[3]
    $ws.Cells.Item($rowNum, 12).Value = $r[4]
    $ws.Cells.Item($rowNum, 13).Value = $r[5]
    $ws.Cells.Item($rowNum, 14).Value = "$/docena de atados"
    $ws.Cells.Item($rowNum, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($rowNum, 16).Value = $r[6]
    $ws.Cells.Item($rowNum, 17).Value = 3
    $ws.Cells.Item($rowNum, 18).Value = "Hortaliza"
}
